$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task detail")

# --- Update C10 FIRST: append " - ok" to the "Hiện avatar, hiện thời gian" line ---
$c10 = $ws.Range("C10").Value2
$c10New = $c10.Replace("Hiện avatar, hiện thời gian`n", "Hiện avatar, hiện thời gian - ok`n")
$ws.Range("C10").Value = $c10New

# --- Update C7 SECOND: "Load list user" task detail text gets " - ok" appended to each bullet ---
$ws.Range("C7").Value = "Load list user`n + Show name + avatar - ok`n + Load thêm khi cuộc xuống cuối container - ok"

# --- Update the sheet view: frozen-pane scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select()
